# Generate Report for Handoff
# Refresh the "Ready for handoff" rows (0907a553, 1a276623, 3e4bf032,
# 927ff260, cd133b31, e6ead3b7 -> sheet rows 7,8,9,11,13,14) with the
# newly generated handoff timestamps and priority/handoff type.

$wb = $excel.ActiveWorkbook

$rows = @(7, 8, 9, 11, 13, 14)

# Overview sheet: column G = "Latest HO Xliff Generate Date"
$overview = $wb.Worksheets.Item("Overview")
foreach ($r in $rows) {
    $overview.Cells.Item($r, 7).Value = "2016-08-22 11:10:59"
}

# zh-cn sheet: column E = "Priority", column H = "Latest Handoff Datetime"
$zhcn = $wb.Worksheets.Item("zh-cn")
foreach ($r in $rows) {
    $zhcn.Cells.Item($r, 5).Value = "ht"
    $zhcn.Cells.Item($r, 8).Value = "2016-08-22 11:10:53"
}

# de-de sheet: column E = "Priority", column H = "Latest Handoff Datetime"
$dede = $wb.Worksheets.Item("de-de")
foreach ($r in $rows) {
    $dede.Cells.Item($r, 5).Value = "ht"
    $dede.Cells.Item($r, 8).Value = "2016-08-22 11:10:59"
}
